$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I (shifts old I..L to J..M).
$ws.Range("I1").EntireColumn.Insert()
$ws.Columns("I:I").ColumnWidth = $ws.Columns("H:H").ColumnWidth

# Move "Generated Date : " label from H5 to the newly inserted I5 (same style as H5).
$ws.Range("I5").Value2 = $ws.Range("H5").Value2
$ws.Range("I5").Style = $ws.Range("H5").Style
$ws.Range("H5").Value2 = $null

# Move "Generated By :" label from H6 to the newly inserted I6 (same style as H6).
$ws.Range("I6").Value2 = $ws.Range("H6").Value2
$ws.Range("I6").Style = $ws.Range("H6").Style
$ws.Range("H6").Value2 = $null

# New "Created By" name-entry field next to "Generated By :".
$ws.Range("J6").Value2 = $null
$ws.Range("J6").Font.Name = "Times New Roman"
$ws.Range("J6").Font.Size = 12
$ws.Range("J6").HorizontalAlignment = -4131
$ws.Range("J6").VerticalAlignment = -4108

# New "Created By" header label in the table header row.
$ws.Range("I8").Value2 = "Created By"
$ws.Range("I8").Style = $ws.Range("J8").Style

# Expand the title merge to include the newly inserted column.
$ws.Range("B2:I2").UnMerge()
$ws.Range("B2:J2").Merge()

$ws.Range("J6").Select()
